# Daily attendance processing - 2025-12-22 21:51:43
# Rotate the comma-separated "Recorded By" values in column G so that the
# last entry in the list moves to the front (right rotation by one).
# Cells with only a single value are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Row + $usedRange.Rows.Count - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = "Recorded By"
    $val = $cell.Value2

    if ($null -ne $val -and $val -ne "") {
        $parts = $val -split ", "
        if ($parts.Count -gt 1) {
            $rotated = @($parts[-1]) + $parts[0..($parts.Count - 2)]
            $cell.Value = [string]::Join(", ", $rotated)
        }
    }
}
